# Apply the commit's changes:
#  - Update B1:B24 on the "SoCtMbCtbDP" sheet to 0.6 (B1 becomes numeric, was a text label)
#  - Make "SoCtMbCtbDP" the active/selected sheet (instead of "About")
#  - Update selection on "SoCtMbCtbDP" sheet to B1:B24 (active cell B1)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("SoCtMbCtbDP")

# Update the data column values: B1 (previously a text label) and B2:B24 (previously 0.5)
# all become the calibrated value 0.6.
$wsData.Range("B1:B24").Value = 0.6

# Make the data sheet the active (selected) tab instead of "About",
# then select B1:B24 so B1 (not B2) is the active cell.
$wsData.Activate()
$wsData.Range("B1:B24").Select() | Out-Null

$wb.Save()
